# Applies the "Updated cryptos list" GitHub Actions refresh: for each
# affected row, rewrites the Coin / Link / Price / Volume(1h) cells with
# the freshly scraped values. Rows 14/15 and 42/43 also swap which coin
# occupies which rank, so Coin (B) and Link (C) are rewritten there too.
#
# Price (D) values are written as text (quote-prefixed display, e.g.
# "1.00" / "0.999") so Excel doesn't reinterpret them as numbers and
# strip the trailing/grouping formatting the source sheet relies on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="75.626.40"; E="  +8.69%  " },
    @{ Row=3; D="2.735.20"; E="  +12.40%  " },
    @{ Row=4; D="1.00" },
    @{ Row=5; D="188.02"; E="  +12.61%  " },
    @{ Row=6; D="592.71"; E="  +4.75%  " },
    @{ Row=7; E="  -0.09%  " },
    @{ Row=8; D="0.544"; E="  +5.80%  " },
    @{ Row=9; D="0.197"; E="  +15.47%  " },
    @{ Row=10; D="2.734.07"; E="  +12.44%  " },
    @{ Row=11; E="  +1.40%  " },
    @{ Row=12; D="0.363"; E="  +8.36%  " },
    @{ Row=13; E="  +2.39%  " },
    @{ Row=14; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="3.183.36"; E="  +10.43%  " },
    @{ Row=15; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="75.473.59"; E="  +8.60%  " },
    @{ Row=16; E="  +7.48%  " },
    @{ Row=17; D="27.15"; E="  +13.01%  " },
    @{ Row=18; D="2.712.71"; E="  +11.34%  " },
    @{ Row=19; D="9.54"; E="  +33.15%  " },
    @{ Row=20; D="12.26"; E="  +12.93%  " },
    @{ Row=21; D="378.53"; E="  +10.05%  " },
    @{ Row=22; D="2.32"; E="  +16.92%  " },
    @{ Row=23; D="4.11"; E="  +6.21%  " },
    @{ Row=24; E="  +5.02%  " },
    @{ Row=25; D="71.34"; E="  +7.80%  " },
    @{ Row=26; D="1.00"; E="  -0.04%  " },
    @{ Row=27; E="  +11.56%  " },
    @{ Row=28; D="9.67"; E="  +13.74%  " },
    @{ Row=29; D="2.845.44"; E="  +11.31%  " },
    @{ Row=30; E="  -1.95%  " },
    @{ Row=31; D="0.0₃0992"; E="  +16.22%  " },
    @{ Row=32; D="526.99"; E="  +15.82%  " },
    @{ Row=33; D="1.42"; E="  +13.94%  " },
    @{ Row=34; D="7.92"; E="  +7.15%  " },
    @{ Row=35; E="  +10.20%  " },
    @{ Row=36; E="  -0.07%  " },
    @{ Row=37; D="0.120"; E="  +7.81%  " },
    @{ Row=38; D="161.39"; E="  +1.78%  " },
    @{ Row=39; D="19.62"; E="  +7.44%  " },
    @{ Row=40; D="19.38"; E="  +1.36%  " },
    @{ Row=41; E="  +0.05%  " },
    @{ Row=42; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="174.11"; E="  +27.84%  " },
    @{ Row=43; B="RenderToken"; C="https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"; D="5.10"; E="  +15.25%  " },
    @{ Row=44; D="1.73"; E="  +13.39%  " },
    @{ Row=45; D="0.335"; E="  +10.25%  " },
    @{ Row=46; E="  +12.96%  " },
    @{ Row=47; D="2.43"; E="  +15.89%  " },
    @{ Row=48; D="39.16"; E="  +3.22%  " },
    @{ Row=49; D="0.0853"; E="  +18.23%  " },
    @{ Row=50; E="  +9.83%  " },
    @{ Row=51; D="0.553"; E="  +12.63%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey("B")) {
        $ws.Cells.Item($row, 2).Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Cells.Item($row, 3).Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        # Force text format first - several prices ("1.00", "0.999", ...)
        # would otherwise be auto-parsed as numbers by Excel.
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}

Write-Host ("Updated {0} rows" -f $updates.Count)
